$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9219396710395813
$ws.Range("B1").Value = 1.114791989326477
$ws.Range("C1").Value = 8.942621231079102
$ws.Range("D1").Value = 2.29966402053833
$ws.Range("E1").Value = 1.25879955291748
